# Update the "Campaign Benchmarks - Historic Context" table on slide 2
# and the footnote/speaker-notes paragraph on slide 3.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2: table of historic campaigns
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tbl = $s2.Shapes.Item(2).Table

# Header row
$tbl.Cell(1, 4).Shape.TextFrame.TextRange.Text = 'Media Spend Insight'
$tbl.Cell(1, 5).Shape.TextFrame.TextRange.Text = 'Conversion Highlight'
$tbl.Cell(1, 6).Shape.TextFrame.TextRange.Text = 'Reference'

# Row 2: Ready Set Ford
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = 'Not yet disclosed; leaning on storytelling + earned media'
$tbl.Cell(2, 5).Shape.TextFrame.TextRange.Text = 'Google Trends first measurable index (1) and 24.5K YouTube anthem views'
$tbl.Cell(2, 6).Shape.TextFrame.TextRange.Text = 'notes/FordNews_Introducing_Ready_Set_Ford.txt:19, data/external/youtube_metrics.csv'

# Row 3: From America, For America
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = 'Employee-pricing reassurance during tariff pressure'
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = 'Undisclosed; cited as Ford’s biggest U.S. push YTD'
$tbl.Cell(3, 5).Shape.TextFrame.TextRange.Text = 'Nationwide employee pricing CRM blast boosted lead volume'
$tbl.Cell(3, 6).Shape.TextFrame.TextRange.Text = 'notes/USA_Today_Ready_Set_Ford.txt:68-76'

# Row 4: Built Ford Proud
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = 'Craftsmanship/heritage storytelling for truck/SUV loyalists'
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Text = 'Est. ~$100M U.S. spend across TV/digital (AdAge)'
$tbl.Cell(4, 5).Shape.TextFrame.TextRange.Text = 'Lifted F-Series loyalty scores; ongoing dealer activations'

# Row 5: Go Further
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = 'Prior global brand platform baseline'
$tbl.Cell(5, 4).Shape.TextFrame.TextRange.Text = 'Global spend not released; served as brand reset'
$tbl.Cell(5, 5).Shape.TextFrame.TextRange.Text = 'Provided recall benchmark still referenced in 2025 coverage'
$tbl.Cell(5, 6).Shape.TextFrame.TextRange.Text = 'notes/USA_Today_Ready_Set_Ford.txt:68-72'

# Row 6: Swap Your Ride
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = 'Comparison drive & testimonial program'
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = 'Broadcast-heavy spend ~ $80M (AdAge archives)'
$tbl.Cell(6, 5).Shape.TextFrame.TextRange.Text = 'Generated testimonial conversions; hashtag revived (#SwapYourRide)'
$tbl.Cell(6, 6).Shape.TextFrame.TextRange.Text = 'notes/research_notes.md:61, external_trade_press:Swap_Your_Ride_2007'

# ---------------------------------------------------------------------
# Slide 3: replace the italic "See <file> for detailed notes..." line
# with a bold "Speaker Notes:" lead-in followed by plain body copy.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr = $s3.Shapes.Item(1).TextFrame.TextRange

# Clear the existing (italic, partly Courier) run(s) first and drop the
# italic override while the range is empty, then type the replacement
# text into the now-plain range so it does not inherit the old runs'
# formatting.
$tr.Text = ''
$tr.Font.Italic = $false
$tr.Text = 'Speaker Notes: - Use this table to show Ready Set Ford’s early traction versus legacy campaigns; highlight how we will track earned reach (Trends, YouTube) now while older programs leaned on paid volumes. - Call out where spend data is undisclosed and how the new BI dashboard will backfill with media + conversion KPIs as RXF matures. - Emphasise that historic campaigns serve as baseline scenarios for loyalty, pricing elasticity, and testimonial-driven conversions.'

$lead = $tr.Characters(1, 14)
$lead.Font.Bold = $true
